$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.861.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.463.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.36%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.33%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.462.35"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.38%  "
$ws.Range("E10").Value = "  -3.74%  "
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.88%  "
$ws.Range("E13").Value = "  -5.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.913.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.826.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("E17").Value = "  -5.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.442.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.90%  "
$ws.Range("E28").Value = "  -49.86%  "
$ws.Range("E30").Value = "  -7.97%  "
$ws.Range("E31").Value = "  -6.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.54%  "
$ws.Range("E34").Value = "  -8.17%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("E37").Value = "  -12.43%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("E39").Value = "  -5.89%  "
$ws.Range("E40").Value = "  -8.98%  "
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.02%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.71%  "
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.324"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.78%  "
$ws.Range("E45").Value = "  -7.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "140.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.77%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.32%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.508"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0252"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -10.54%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0730"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.15%  "

Write-Host "Applied all changes"
